$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.663.00"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "2.503.95"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  +0.03%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "574.28"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.75%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "166.40"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("D9").Value = "2.503.45"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("E11").Value = "  -0.08%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.359"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +3.05%  "
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "2.959.40"
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("D15").Value = "69.464.82"
$ws.Range("E15").Value = "  -0.79%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.0000177"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.39%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "24.71"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -1.84%  "
$ws.Range("D18").Value = "2.499.22"
$ws.Range("E18").Value = "  -1.01%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "11.22"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("E20").Value = "  -3.61%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "348.97"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.71%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "3.91"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("E24").Value = "  -0.05%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "70.82"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +2.28%  "
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("E27").Value = "  -3.33%  "
$ws.Range("D28").Value = "2.630.77"
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").Value = "0.0₃0891"
$ws.Range("E30").Value = "  -2.17%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "7.87"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.67%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "457.93"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -1.72%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.20"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -6.19%  "
$ws.Range("E34").Value = "  -1.68%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "157.26"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.116"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -3.71%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "19.05"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("E39").Value = "  -0.80%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("E42").Value = "  -2.43%  "
$ws.Range("E43").Value = "  -0.18%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "38.16"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("E45").Value = "  -5.39%  "
$ws.Range("E46").Value = "  -7.91%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "141.12"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("E49").Value = "  -2.81%  "
$ws.Range("E50").Value = "  -0.55%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.578"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.88%  "
